$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = "Roblox's DevForum, Discord, Reddit, LinkedIn"
$ws.Range("D13").Value = "Roblox's DevForum"
$ws.Range("E13").Value = "Roblox's DevForum"
$ws.Range("F13").Value = "Roblox's DevForum"
$ws.Range("G13").Value = "Roblox's DevForum"

$ws.Range("G14").Select()
